$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.195025043834562
$ws.Range("C2").Value = 0.1929351714631196
$ws.Range("D2").Value = 0.07675989385728599
$ws.Range("E2").Value = 0.07729536982203911
$ws.Range("G2").Value = 0.002572860717651765
$ws.Range("K2").Value = 1.309855205363078
$ws.Range("M2").Value = 0.4099677081541913
$ws.Range("N2").Value = 4.324869992211319
$ws.Range("B3").Value = 1.133608488428848
$ws.Range("C3").Value = 0.1798040361767903
$ws.Range("D3").Value = 0.06987593982898943
$ws.Range("E3").Value = 0.06974987679675593
$ws.Range("G3").Value = 0.002578550013180779
$ws.Range("K3").Value = 1.237781535250832
$ws.Range("M3").Value = 0.3819759296136596
$ws.Range("N3").Value = 4.238471299622972
$ws.Range("B4").Value = 1.096674237852653
$ws.Range("C4").Value = 0.1718682979771131
$ws.Range("D4").Value = 0.06569231346102811
$ws.Range("E4").Value = 0.06516241655330646
$ws.Range("G4").Value = 0.002582222482562702
$ws.Range("K4").Value = 1.194384315689092
$ws.Range("M4").Value = 0.3650564059615036
$ws.Range("N4").Value = 4.185769544798717
$ws.Range("B5").Value = 1.081816758024416
$ws.Range("C5").Value = 0.1686658855076359
$ws.Range("D5").Value = 0.06399812098496227
$ws.Range("E5").Value = 0.06330413337100538
$ws.Range("G5").Value = 0.002583764274657237
$ws.Range("K5").Value = 1.176912977158366
$ws.Range("M5").Value = 0.358227950883041
$ws.Range("N5").Value = 4.164378877138091
$ws.Range("B6").Value = 1.079361328800985
$ws.Range("C6").Value = 0.168136015406759
$ws.Range("D6").Value = 0.06371744077091535
$ws.Range("E6").Value = 0.06299623129783072
$ws.Range("G6").Value = 0.00258402302480603
$ws.Range("K6").Value = 1.174024706223435
$ws.Range("M6").Value = 0.3570980751185644
$ws.Range("N6").Value = 4.160832103315272
$ws.Range("B7").Value = 1.096473083127393
$ws.Range("C7").Value = 0.1718249822804978
$ws.Range("D7").Value = 0.06566942202364601
$ws.Range("E7").Value = 0.06513731036296377
$ws.Range("G7").Value = 0.002582243092037492
$ws.Range("K7").Value = 1.194147829686301
$ws.Range("M7").Value = 0.3649640473853069
$ws.Range("N7").Value = 4.185480717883223
$ws.Range("B8").Value = 1.173686647276185
$ws.Range("C8").Value = 0.1883809686771656
$ws.Range("D8").Value = 0.07437723914951277
$ws.Range("E8").Value = 0.07468405716203108
$ws.Range("G8").Value = 0.002574785289764738
$ws.Range("K8").Value = 1.284825273961019
$ws.Range("M8").Value = 0.4002600680944752
$ws.Range("N8").Value = 4.295006611737534
$ws.Range("B9").Value = 1.331337249307239
$ws.Range("C9").Value = 0.2218748531168444
$ws.Range("D9").Value = 0.09180474945503647
$ws.Range("E9").Value = 0.09378072326197895
$ws.Range("G9").Value = 0.002561575007495277
$ws.Range("K9").Value = 1.469537125329623
$ws.Range("M9").Value = 0.4716425372866979
$ws.Range("N9").Value = 4.512627668471993
$ws.Range("B10").Value = 1.451089078158304
$ws.Range("C10").Value = 0.2471416672106272
$ws.Range("D10").Value = 0.1048370893344384
$ws.Range("E10").Value = 0.1080615253800659
$ws.Range("G10").Value = 0.002552721060895196
$ws.Range("K10").Value = 1.609601663125773
$ws.Range("M10").Value = 0.5254757577652072
$ws.Range("N10").Value = 4.674379392255418
$ws.Range("B11").Value = 1.50644806759982
$ws.Range("C11").Value = 0.2587865287289901
$ws.Range("D11").Value = 0.1108185801734436
$ws.Range("E11").Value = 0.1146174215653204
$ws.Range("G11").Value = 0.002548875820673063
$ws.Range("K11").Value = 1.674301519117876
$ws.Range("M11").Value = 0.550282519585906
$ws.Range("N11").Value = 4.748399490264546
$ws.Range("B12").Value = 1.527540161373508
$ws.Range("C12").Value = 0.2632183751721584
$ws.Range("D12").Value = 0.113091471203191
$ws.Range("E12").Value = 0.1171088906032551
$ws.Range("G12").Value = 0.002547445791690182
$ws.Range("K12").Value = 1.698945708398639
$ws.Range("M12").Value = 0.5597230227212009
$ws.Range("N12").Value = 4.77649412679466
$ws.Range("B13").Value = 1.522991841556291
$ws.Range("C13").Value = 0.2622629028273877
$ws.Range("D13").Value = 0.1126016118671913
$ws.Range("E13").Value = 0.1165719071712772
$ws.Range("G13").Value = 0.00254775261674617
$ws.Range("K13").Value = 1.693631710658678
$ws.Range("M13").Value = 0.557687746978047
$ws.Range("N13").Value = 4.77044053201854
$ws.Range("B14").Value = 1.508180733954475
$ws.Range("C14").Value = 0.2591506918023185
$ws.Range("D14").Value = 0.1110054141578871
$ws.Range("E14").Value = 0.1148222160221053
$ws.Range("G14").Value = 0.002548757649676705
$ws.Range("K14").Value = 1.676326116994005
$ws.Range("M14").Value = 0.551058253152874
$ws.Range("N14").Value = 4.750709543283335
$ws.Range("B15").Value = 1.499125337449414
$ws.Range("C15").Value = 0.2572472773889842
$ws.Range("D15").Value = 0.1100287229066481
$ws.Range("E15").Value = 0.1137516480211787
$ws.Range("G15").Value = 0.002549376653000695
$ws.Range("K15").Value = 1.665744731480231
$ws.Range("M15").Value = 0.5470036120078703
$ws.Range("N15").Value = 4.738632239382582
$ws.Range("B16").Value = 1.44748916752178
$ws.Range("C16").Value = 0.2463837260615662
$ws.Range("D16").Value = 0.104447269224849
$ws.Range("E16").Value = 0.1076343080597013
$ws.Range("G16").Value = 0.002552976014497907
$ws.Range("K16").Value = 1.605393366481849
$ws.Range("M16").Value = 0.5238610590687642
$ws.Range("N16").Value = 4.669550979256115
$ws.Range("B17").Value = 1.416039569857901
$ws.Range("C17").Value = 0.2397582786417445
$ws.Range("D17").Value = 0.1010369524289416
$ws.Range("E17").Value = 0.1038969939622731
$ws.Range("G17").Value = 0.002555230727000985
$ws.Range("K17").Value = 1.568623346468485
$ws.Range("M17").Value = 0.50974596335206
$ws.Range("N17").Value = 4.627285419585405
$ws.Range("B18").Value = 1.398033619885155
$ws.Range("C18").Value = 0.2359616660262134
$ws.Range("D18").Value = 0.09908040574973143
$ws.Range("E18").Value = 0.1017529683899099
$ws.Range("G18").Value = 0.002556544762002845
$ws.Range("K18").Value = 1.547566663188377
$ws.Range("M18").Value = 0.5016571801080687
$ws.Range("N18").Value = 4.603016546276109
$ws.Range("B19").Value = 1.391951319387829
$ws.Range("C19").Value = 0.2346786168805579
$ws.Range("D19").Value = 0.09841880128205105
$ws.Range("E19").Value = 0.1010279860649987
$ws.Range("G19").Value = 0.002556992627747983
$ws.Range("K19").Value = 1.540453043086359
$ws.Range("M19").Value = 0.4989235522555333
$ws.Range("N19").Value = 4.59480653773096
$ws.Range("B20").Value = 1.419378827122671
$ws.Range("C20").Value = 0.24046209877352
$ws.Range("D20").Value = 0.1013994698502501
$ws.Range("E20").Value = 0.1042942576194932
$ws.Range("G20").Value = 0.00255498893166284
$ws.Range("K20").Value = 1.572527992209928
$ws.Range("M20").Value = 0.511245443008491
$ws.Range("N20").Value = 4.631780387925204
$ws.Range("B21").Value = 1.512527603652586
$ws.Range("C21").Value = 0.2600642173845245
$ws.Range("D21").Value = 0.1114740419392461
$ws.Range("E21").Value = 0.1153358986302138
$ws.Range("G21").Value = 0.0025484617405584
$ws.Range("K21").Value = 1.681405269314723
$ws.Range("M21").Value = 0.5530042199196572
$ws.Range("N21").Value = 4.75650323908917
$ws.Range("B22").Value = 1.574157604680011
$ws.Range("C22").Value = 0.2730049131729402
$ws.Range("D22").Value = 0.1181040966072828
$ws.Range("E22").Value = 0.1226042484513954
$ws.Range("G22").Value = 0.00254434778446964
$ws.Range("K22").Value = 1.753401892071793
$ws.Range("M22").Value = 0.5805688097208304
$ws.Range("N22").Value = 4.838395567377006
$ws.Range("B23").Value = 1.541195099833772
$ws.Range("C23").Value = 0.2660861973834869
$ws.Range("D23").Value = 0.1145612571789485
$ws.Range("E23").Value = 0.1187201219250653
$ws.Range("G23").Value = 0.002546529629840795
$ws.Range("K23").Value = 1.714898399508911
$ws.Range("M23").Value = 0.5658317623568934
$ws.Range("N23").Value = 4.794652840032825
$ws.Range("B24").Value = 1.417868915932672
$ws.Range("C24").Value = 0.2401438629648567
$ws.Range("D24").Value = 0.1012355630007278
$ws.Range("E24").Value = 0.1041146403987909
$ws.Range("G24").Value = 0.002555098192276937
$ws.Range("K24").Value = 1.570762443855415
$ws.Range("M24").Value = 0.5105674468091905
$ws.Range("N24").Value = 4.629748118931843
$ws.Range("B25").Value = 1.288008304897801
$ws.Range("C25").Value = 0.2127005379638263
$ws.Range("D25").Value = 0.08705107369578968
$ws.Range("E25").Value = 0.08857215693641507
$ws.Range("G25").Value = 0.002564998413873143
$ws.Range("K25").Value = 1.418813884033796
$ws.Range("M25").Value = 0.4520927824124783
$ws.Range("N25").Value = 4.453437516539537
